# Auto-generated edit script: refresh the crypto price / 1h-volume table
# (Bitcoin .. the 50th coin) with newly scraped values.
#
# Columns: B=Coin name, C=Link, D=Price (text), E=Volume(1h) (text)
#
# Price cells whose text looks like a plain decimal number (e.g. "227.10")
# would otherwise be auto-converted to a real number by Excel (losing the
# trailing zero / exact text). A leading apostrophe forces Excel to keep
# storing them as text, exactly like the original data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.797.91'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '2.101.53'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''227.10'
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = '''0.617'
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("D7").Value = '''62.04'
$ws.Range("E7").Value = '  +2.99%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.388'
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("D10").Value = '''0.0839'
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '''15.83'
$ws.Range("E12").Value = '  +6.58%  '
$ws.Range("D13").Value = '2.412.23'
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").Value = '''21.96'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").Value = '''0.801'
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").Value = '''5.49'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D17").Value = '2.101.27'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '38.791.28'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").Value = '''71.52'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").Value = '''6.06'
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").Value = '0.0₃0845'
$ws.Range("E21").Value = '  +1.65%  '
$ws.Range("D22").Value = '''227.33'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '''2.41'
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("D25").Value = '''2.31'
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''170.81'
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''9.51'
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").Value = '''0.136'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '''1.42'
$ws.Range("E29").Value = '  +3.67%  '
$ws.Range("D30").Value = '''19.29'
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("D31").Value = '''2.54'
$ws.Range("E31").Value = '  +9.13%  '
$ws.Range("D32").Value = '''0.121'
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").Value = '''7.24'
$ws.Range("E33").Value = '  +14.22%  '
$ws.Range("D34").Value = '''4.57'
$ws.Range("E34").Value = '  +1.54%  '
$ws.Range("D35").Value = '''4.75'
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").Value = '''0.0614'
$ws.Range("E36").Value = '  +1.44%  '
$ws.Range("D37").Value = '''2.37'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '''3.50'
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '''18.01'
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("D41").Value = '''0.0227'
$ws.Range("E41").Value = '  +3.28%  '
$ws.Range("D42").Value = '''101.51'
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("D43").Value = '1.524.11'
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("E44").Value = '  +7.45%  '
$ws.Range("D45").Value = '''2.80'
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = '''7.78'
$ws.Range("E46").Value = '  +1.23%  '
$ws.Range("D47").Value = '''0.0914'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("E48").Value = '  +5.39%  '
$ws.Range("D49").Value = '''4.16'
$ws.Range("E49").Value = '  +2.17%  '
$ws.Range("D50").Value = '''2.96'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").Value = '2.299.71'
$ws.Range("E51").Value = '  +0.98%  '
